$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Range("H6").Value = 210.90909
$ws.Range("I6").Value = 217
$ws.Range("K6").Value = 651
$ws.Range("M6").Value = -539
# Row 12
$ws.Range("H12").Value = 1
$ws.Range("J12").Value = 1
$ws.Range("L12").Value = 1
$ws.Range("N12").Value = -341
# Row 28
$ws.Range("H28").Value = 138.33333
$ws.Range("I28").Value = 138.33333
$ws.Range("K28").Value = 138.33333
$ws.Range("M28").Value = 346.66667
# Row 68
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").ClearContents()
# Row 71
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").ClearContents()
# Row 116
$ws.Range("H116").Value = 3498.3333
$ws.Range("I116").Value = 2747.5
$ws.Range("K116").Value = 2747.5
$ws.Range("M116").Value = 694.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 28
$ws.Range("H28").Value = 9000
$ws.Range("I28").Value = 9000
$ws.Range("K28").Value = 9000
$ws.Range("M28").Value = -8808
# Row 55
$ws.Range("H55").Value = 26000
$ws.Range("J55").Value = 26000
$ws.Range("L55").Value = 26000
$ws.Range("N55").Value = -26630
# Row 99
$ws.Range("H99").Value = 9000
$ws.Range("I99").Value = 9000
$ws.Range("K99").Value = 9000
$ws.Range("M99").Value = -6005
# Row 110
$ws.Range("H110").Value = 3179.6924
$ws.Range("I110").Value = 2857.5
$ws.Range("K110").Value = 2857.5
$ws.Range("M110").Value = -812.5
# Row 119
$ws.Range("H119").Value = 65999.60000000001
$ws.Range("J119").Value = 65999.60000000001
$ws.Range("L119").Value = 65999.60000000001
$ws.Range("N119").Value = -75675.60000000001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 107
$ws.Range("H107").Value = 1094
$ws.Range("I107").Value = 1094
$ws.Range("K107").Value = 1094
$ws.Range("M107").Value = 826

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 420.2
$ws.Range("I22").Value = 467
$ws.Range("J22").Value = 350
$ws.Range("K22").Value = 467
$ws.Range("L22").Value = 350
$ws.Range("M22").Value = -117
$ws.Range("N22").Value = -1050
# Row 58
$ws.Range("H58").Value = 5162.375
$ws.Range("I58").Value = 4400
$ws.Range("J58").Value = 5416.5
$ws.Range("K58").Value = 4400
$ws.Range("L58").Value = 5416.5
$ws.Range("M58").Value = -4197
$ws.Range("N58").Value = -5822.5
# Row 133
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
# Row 134
$ws.Range("H134").Value = 1204.25
$ws.Range("I134").Value = 1162
$ws.Range("K134").Value = 3486
$ws.Range("M134").Value = -951
# Row 136
$ws.Range("H136").Value = 5162.375
$ws.Range("I136").Value = 4400
$ws.Range("J136").Value = 5416.5
$ws.Range("K136").Value = 13200
$ws.Range("L136").Value = 16249.5
$ws.Range("M136").Value = -10650
$ws.Range("N136").Value = -21349.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 34
$ws.Range("H34").Value = 966.6667
$ws.Range("I34").Value = 500
$ws.Range("J34").Value = 1060
$ws.Range("K34").Value = 1500
$ws.Range("L34").Value = 3180
$ws.Range("M34").Value = -1416
$ws.Range("N34").Value = -3348
# Row 38
$ws.Range("H38").Value = 133.25
$ws.Range("I38").Value = 63.8
$ws.Range("J38").Value = 249
$ws.Range("K38").Value = 191.4
$ws.Range("L38").Value = 747
$ws.Range("M38").Value = 155.6
$ws.Range("N38").Value = -1441
# Row 39
$ws.Range("H39").Value = 5331.6665
$ws.Range("J39").Value = 5998
$ws.Range("L39").Value = 17994
$ws.Range("N39").Value = -18582
# Row 54
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()
# Row 55
$ws.Range("H55").Value = 692.5
$ws.Range("J55").Value = 1000
$ws.Range("L55").Value = 3000
$ws.Range("N55").Value = -3354
# Row 58
$ws.Range("H58").Value = 2290
# Row 123
$ws.Range("H123").Value = 4597
$ws.Range("I123").Value = 4996.25
$ws.Range("K123").Value = 14988.75
$ws.Range("M123").Value = -12538.75
# Row 139
$ws.Range("H139").Value = 3798.4
$ws.Range("I139").Value = 2332.6667
$ws.Range("J139").Value = 5997
$ws.Range("K139").Value = 6998.000100000001
$ws.Range("L139").Value = 17991
$ws.Range("M139").Value = -1858.000100000001
$ws.Range("N139").Value = -28271

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 557.5454999999999
$ws.Range("I2").Value = 769.7143
$ws.Range("J2").Value = 186.25
$ws.Range("K2").Value = 769.7143
$ws.Range("L2").Value = 186.25
$ws.Range("M2").Value = -656.7143
$ws.Range("N2").Value = -412.25
# Row 80
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("N80").ClearContents()
# Row 83
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("N83").ClearContents()
# Row 113
$ws.Range("H113").Value = 962
$ws.Range("I113").Value = 962
$ws.Range("K113").Value = 962
$ws.Range("M113").Value = 1208
# Row 122
$ws.Range("H122").Value = 1999
$ws.Range("I122").Value = 1999
$ws.Range("K122").Value = 5997
$ws.Range("M122").Value = -3547
# Row 123
$ws.Range("H123").Value = 24225
$ws.Range("J123").Value = 24225
$ws.Range("L123").Value = 24225
$ws.Range("N123").Value = -29125

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 82
$ws.Range("H82").Value = 0
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("M82").ClearContents()
$ws.Range("N82").ClearContents()
# Row 85
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("M85").ClearContents()
$ws.Range("N85").ClearContents()
# Row 100
$ws.Range("H100").Value = 4299.6665
$ws.Range("I100").Value = 4299.6665
$ws.Range("K100").Value = 4299.6665
$ws.Range("M100").Value = -3758.6665

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 16
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()
# Row 70
$ws.Range("H70").Value = 35000
$ws.Range("J70").Value = 35000
$ws.Range("L70").Value = 35000
$ws.Range("N70").Value = -35630
# Row 73
$ws.Range("H73").Value = 35000
$ws.Range("J73").Value = 35000
$ws.Range("L73").Value = 35000
$ws.Range("N73").Value = -37184
# Row 81
$ws.Range("H81").Value = 1900
$ws.Range("J81").Value = 1900
$ws.Range("L81").Value = 3800
$ws.Range("N81").Value = -5922
# Row 84
$ws.Range("H84").Value = 1900
$ws.Range("J84").Value = 1900
$ws.Range("L84").Value = 19000
$ws.Range("N84").Value = -29608
